# Second commit for the presentation
#
# The only user-visible content change in this commit is on slide 5
# (sldId 277), shape id 15 ("Content Placeholder 14"): the second
# paragraph ("When we have design ideas, we'll show them to you right
# there. ") was removed, leaving only the first sentence ("Open the
# Design Ideas pane for instant slide makeovers. ") in a single
# paragraph.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(5)
$tr = $shp.TextFrame.TextRange

# Replace the whole text with just the first sentence, then re-append
# the trailing ". " via InsertAfter so the run layout mirrors how
# PowerPoint splits a run after an in-place edit (instead of leaving a
# stray empty trailing paragraph behind).
$tr.Text = "Open the Design Ideas pane for instant slide makeovers"
[void]$tr.InsertAfter(". ")
